# Edit slide 6 ("NSubstitute"): insert a new bullet/paragraph before the
# "Checking calls to properties requires throw-away variables" paragraph,
# inside the "Cons" content placeholder (Shapes.Item(4)).
#
# New paragraph text (as two runs, matching the authored edit):
#   "Extension methods can be called an anything, but throw exceptions if "
#   "not mocks"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(4)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Locate the paragraph that currently reads "Checking calls to properties..."
# (it is the 3rd paragraph in this placeholder: Cons / Extension methods
# invade.../ Checking calls.../ Sometimes hard...).
$targetPara = $tr.Paragraphs(3, 1)

# Insert the full new paragraph text (plus trailing paragraph mark) right
# before it. This new text inherits the run formatting of the paragraph it
# is adjacent to (lang="en-US").
$targetPara.InsertBefore("Extension methods can be called an anything, but throw exceptions if not mocks`r")

# Re-fetch the freshly created paragraph (still at position 3).
$newPara = $tr.Paragraphs(3, 1)

# Split the trailing "not mocks" off into its own run, matching the
# authored two-run paragraph structure.
$splitLen = 9
$splitStart = $newPara.Start + $newPara.Length - $splitLen - 1
$tail = $tr.Characters($splitStart, $splitLen)
$tail.Delete()
$head = $tr.Paragraphs(3, 1)
$head.InsertAfter("not mocks")
